# edit.ps1
# "always rebuild courses, output goes to docs"
#
# Content-level changes applied here:
#   1. Under "Required Software/Online Accounts", the "Slack Client ..." and
#      "Firefox or Chrome web browser" bullets are demoted to be sub-bullets
#      of "Software".
#   2. Under the same heading, the "repl.it" and "AU Ed Tech Slack channel
#      #code" bullets are demoted to be sub-bullets of "Accounts".
#   3. Both tables in the document (the Class Sessions schedule and the
#      Assignments and Grading breakdown) switch from the "TableNormal"
#      table style to the "Table" table style.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1 & 2: demote the four list items to the second outline level so they
# nest under "Software" / "Accounts" respectively.
# ---------------------------------------------------------------------
$targets = @(
    "Slack Client (recommend desktop and mobile clients)",
    "Firefox or Chrome web browser",
    "repl.it",
    "AU Ed Tech Slack channel #code"
)

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.Trim()
    if ($targets -contains $text) {
        $p.Range.ListFormat.ListLevelNumber = 2
    }
}

# ---------------------------------------------------------------------
# 3: retarget both tables onto the "Table" table style.
# ---------------------------------------------------------------------
foreach ($t in $d.Tables) {
    $t.Style = "Table"
}

Write-Output "done"
